$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.011.38"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.898.65"
$ws.Range("E3").Value = "  +6.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.02"
$ws.Range("E5").Value = "  +4.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.57"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.35%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "2.897.11"
$ws.Range("E10").Value = "  +6.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("E11").Value = "  +11.03%  "
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.91"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "3.426.79"
$ws.Range("E14").Value = "  +6.72%  "
$ws.Range("D15").Value = "75.890.31"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000192"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.45"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "2.891.79"
$ws.Range("E18").Value = "  +6.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.93"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  +3.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.96"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.17"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.54"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "3.032.93"
$ws.Range("E26").Value = "  +6.13%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000109"
$ws.Range("E29").Value = "  +10.00%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "508.29"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.28"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.24"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "183.35"
$ws.Range("E40").Value = "  +5.93%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.02"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("E45").Value = "  +6.96%  "
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.40"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.582"
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.670"
$ws.Range("E50").Value = "  +11.79%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.77"
$ws.Range("E51").Value = "  +1.40%  "
